$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegistrationUser")

$ws.Range("A6").Value = "RegistrationWithoutFullName"
$ws.Range("B6").Value = "email@abv.bg"
$ws.Range("D6").Value = 123
$ws.Range("E6").Value = 123
